$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Move "Tipo" header from D1 to E1 (new column), carrying the header style,
# then insert the new "MAE" header at D1 with the same style.
$ws.Range("C1").Copy()
$ws.Range("E1").PasteSpecial(-4122)
$ws.Range("E1").Value = "Tipo"

$ws.Range("C1").Copy()
$ws.Range("D1").PasteSpecial(-4122)
$ws.Range("D1").Value = "MAE"

# Move "single" value from D2 to E2, and fill the new MAE metric at D2.
$ws.Range("E2").Value = "single"
$ws.Range("D2").Value = 0.6669627832707453

# Updated MSE / R2 metric values.
$ws.Range("B2").Value = 0.7399047889235683
$ws.Range("C2").Value = 0.9782670581212405
